$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Plg"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("I2").Value = 0.9660838355812051
$ws.Range("J2").Value = 0.9660838355812051
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 28.43067752420701
$ws.Range("R2").Value = 255.8760977178631
$ws.Range("S2").Value = 0.1974301930437226
$ws.Range("T2").Value = 0.1974301930437226

# Row 3
$ws.Range("B3").Value = "Plg"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("I3").Value = 0.9660838355812051
$ws.Range("J3").Value = 0.9660838355812051
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("S3").Value = 0.3438586696277408
$ws.Range("T3").Value = 0.3438586696277409

# Row 4
$ws.Range("B4").Value = "Plg"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 0.9660838355812051
$ws.Range("J4").Value = 0.9660838355812051
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 61.172046192684
$ws.Range("R4").Value = 550.548415734156
$ws.Range("S4").Value = 0.4247949729097417
$ws.Range("T4").Value = 0.4247949729097417

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Plg"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01635133333333333
$ws.Range("H5").Value = 0.049054
$ws.Range("I5").Value = 0.03391616441879487
$ws.Range("J5").Value = 0.03391616441879487
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 0.9981116523584445
$ws.Range("R5").Value = 8.983004871226
$ws.Range("S5").Value = 0.006931153013730826
$ws.Range("T5").Value = 0.006931153013730826

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Plg"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01635133333333333
$ws.Range("H6").Value = 0.049054
$ws.Range("I6").Value = 0.03391616441879487
$ws.Range("J6").Value = 0.03391616441879487
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("Q6").Value = 1.738383271721333
$ws.Range("R6").Value = 15.645449445492
$ws.Range("S6").Value = 0.01207179620069549
$ws.Range("T6").Value = 0.01207179620069549

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Plg"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01635133333333333
$ws.Range("H7").Value = 0.049054
$ws.Range("I7").Value = 0.03391616441879487
$ws.Range("J7").Value = 0.03391616441879487
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 2.147558110479111
$ws.Range("R7").Value = 19.328022994312
$ws.Range("S7").Value = 0.01491321520436855
$ws.Range("T7").Value = 0.01491321520436855
